$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'26.584.15"
$ws.Cells.Item(2, 4).Style = 'Normal'
$ws.Cells.Item(2, 5).Value = '  +0.57%  '

$ws.Cells.Item(3, 4).Value = "'1.840.01"
$ws.Cells.Item(3, 4).Style = 'Normal'
$ws.Cells.Item(3, 5).Value = '  -0.10%  '

$ws.Cells.Item(4, 5).Value = '  -0.03%  '

$ws.Cells.Item(5, 4).Value = "'258.80"
$ws.Cells.Item(5, 4).Style = 'Normal'
$ws.Cells.Item(5, 5).Value = '  -1.11%  '

$ws.Cells.Item(6, 4).Value = "'0.9998"
$ws.Cells.Item(6, 4).Style = 'Normal'
$ws.Cells.Item(6, 5).Value = '  -0.04%  '

$ws.Cells.Item(7, 4).Value = "'0.5269"
$ws.Cells.Item(7, 4).Style = 'Normal'
$ws.Cells.Item(7, 5).Value = '  +1.06%  '

$ws.Cells.Item(8, 4).Value = "'0.3150"
$ws.Cells.Item(8, 4).Style = 'Normal'
$ws.Cells.Item(8, 5).Value = '  -3.30%  '

$ws.Cells.Item(9, 4).Value = "'0.06794"
$ws.Cells.Item(9, 4).Style = 'Normal'
$ws.Cells.Item(9, 5).Value = '  +0.11%  '

$ws.Cells.Item(10, 4).Value = "'18.67"
$ws.Cells.Item(10, 4).Style = 'Normal'
$ws.Cells.Item(10, 5).Value = '  +0.07%  '

$ws.Cells.Item(11, 4).Value = "'0.7785"
$ws.Cells.Item(11, 4).Style = 'Normal'
$ws.Cells.Item(11, 5).Value = '  +0.98%  '

$ws.Cells.Item(12, 5).Value = '  +0.45%  '

$ws.Cells.Item(13, 4).Value = "'1.862.99"
$ws.Cells.Item(13, 4).Style = 'Normal'
$ws.Cells.Item(13, 5).Value = '  +1.28%  '

$ws.Cells.Item(14, 4).Value = "'87.78"
$ws.Cells.Item(14, 4).Style = 'Normal'
$ws.Cells.Item(14, 5).Value = '  -0.06%  '

$ws.Cells.Item(15, 4).Value = "'5.003"
$ws.Cells.Item(15, 4).Style = 'Normal'
$ws.Cells.Item(15, 5).Value = '  +0.04%  '

$ws.Cells.Item(16, 4).Value = "'0.9996"
$ws.Cells.Item(16, 4).Style = 'Normal'
$ws.Cells.Item(16, 5).Value = '  +0.01%  '

$ws.Cells.Item(17, 5).Value = '  -0.67%  '

$ws.Cells.Item(18, 4).Value = "'0.9998"
$ws.Cells.Item(18, 4).Style = 'Normal'
$ws.Cells.Item(18, 5).Value = '  -0.03%  '

$ws.Cells.Item(19, 4).Value = "'0.000007914"
$ws.Cells.Item(19, 4).Style = 'Normal'
$ws.Cells.Item(19, 5).Value = '  -0.62%  '

$ws.Cells.Item(20, 4).Value = "'26.592.24"
$ws.Cells.Item(20, 4).Style = 'Normal'
$ws.Cells.Item(20, 5).Value = '  +0.53%  '

$ws.Cells.Item(21, 4).Value = "'2.070.16"
$ws.Cells.Item(21, 4).Style = 'Normal'
$ws.Cells.Item(21, 5).Value = '  +0.20%  '

$ws.Cells.Item(22, 4).Value = "'4.594"
$ws.Cells.Item(22, 4).Style = 'Normal'
$ws.Cells.Item(22, 5).Value = '  +0.18%  '

$ws.Cells.Item(23, 4).Value = "'5.962"
$ws.Cells.Item(23, 4).Style = 'Normal'
$ws.Cells.Item(23, 5).Value = '  -0.31%  '

$ws.Cells.Item(24, 4).Value = "'9.312"
$ws.Cells.Item(24, 4).Style = 'Normal'
$ws.Cells.Item(24, 5).Value = '  -1.93%  '

$ws.Cells.Item(25, 4).Value = "'143.13"
$ws.Cells.Item(25, 4).Style = 'Normal'
$ws.Cells.Item(25, 5).Value = '  -1.09%  '

$ws.Cells.Item(26, 4).Value = "'2.213"
$ws.Cells.Item(26, 4).Style = 'Normal'
$ws.Cells.Item(26, 5).Value = '  +1.53%  '

$ws.Cells.Item(27, 4).Value = "'1.686"
$ws.Cells.Item(27, 4).Style = 'Normal'
$ws.Cells.Item(27, 5).Value = '  +2.16%  '

$ws.Cells.Item(28, 4).Value = "'16.95"
$ws.Cells.Item(28, 4).Style = 'Normal'
$ws.Cells.Item(28, 5).Value = '  -0.24%  '

$ws.Cells.Item(29, 4).Value = "'110.94"
$ws.Cells.Item(29, 4).Style = 'Normal'
$ws.Cells.Item(29, 5).Value = '  -0.35%  '

$ws.Cells.Item(30, 4).Value = "'4.177"
$ws.Cells.Item(30, 4).Style = 'Normal'
$ws.Cells.Item(30, 5).Value = '  -0.60%  '

$ws.Cells.Item(31, 4).Value = "'0.08717"
$ws.Cells.Item(31, 4).Style = 'Normal'

$ws.Cells.Item(32, 4).Value = "'4.060"
$ws.Cells.Item(32, 4).Style = 'Normal'
$ws.Cells.Item(32, 5).Value = '  -1.59%  '

$ws.Cells.Item(33, 4).Value = "'0.04867"
$ws.Cells.Item(33, 4).Style = 'Normal'
$ws.Cells.Item(33, 5).Value = '  +1.48%  '

$ws.Cells.Item(34, 4).Value = "'0.7307"
$ws.Cells.Item(34, 4).Style = 'Normal'
$ws.Cells.Item(34, 5).Value = '  +1.72%  '

$ws.Cells.Item(35, 4).Value = "'1.137"
$ws.Cells.Item(35, 4).Style = 'Normal'
$ws.Cells.Item(35, 5).Value = '  +0.68%  '

$ws.Cells.Item(36, 4).Value = "'2.859"
$ws.Cells.Item(36, 4).Style = 'Normal'
$ws.Cells.Item(36, 5).Value = '  +0.45%  '

$ws.Cells.Item(37, 4).Value = "'3.086"
$ws.Cells.Item(37, 4).Style = 'Normal'
$ws.Cells.Item(37, 5).Value = '  +0.14%  '

$ws.Cells.Item(38, 4).Value = "'2.242"
$ws.Cells.Item(38, 4).Style = 'Normal'
$ws.Cells.Item(38, 5).Value = '  +1.29%  '

$ws.Cells.Item(39, 4).Value = "'0.01717"
$ws.Cells.Item(39, 4).Style = 'Normal'
$ws.Cells.Item(39, 5).Value = '  -3.58%  '

$ws.Cells.Item(40, 4).Value = "'0.4803"
$ws.Cells.Item(40, 4).Style = 'Normal'
$ws.Cells.Item(40, 5).Value = '  -0.54%  '

$ws.Cells.Item(41, 4).Value = "'0.8942"
$ws.Cells.Item(41, 4).Style = 'Normal'
$ws.Cells.Item(41, 5).Value = '  -0.60%  '

$ws.Cells.Item(42, 4).Value = "'109.76"
$ws.Cells.Item(42, 4).Style = 'Normal'
$ws.Cells.Item(42, 5).Value = '  -2.25%  '

$ws.Cells.Item(43, 4).Value = "'5.919"
$ws.Cells.Item(43, 4).Style = 'Normal'
$ws.Cells.Item(43, 5).Value = '  -2.47%  '

$ws.Cells.Item(44, 4).Value = "'1.0000"
$ws.Cells.Item(44, 4).Style = 'Normal'
$ws.Cells.Item(44, 5).Value = '  -0.01%  '

$ws.Cells.Item(45, 4).Value = "'7.631"
$ws.Cells.Item(45, 4).Style = 'Normal'
$ws.Cells.Item(45, 5).Value = '  -1.15%  '

$ws.Cells.Item(46, 4).Value = "'0.4150"
$ws.Cells.Item(46, 4).Style = 'Normal'
$ws.Cells.Item(46, 5).Value = '  +0.56%  '

$ws.Cells.Item(47, 4).Value = "'8.963"
$ws.Cells.Item(47, 4).Style = 'Normal'
$ws.Cells.Item(47, 5).Value = '  +0.16%  '

$ws.Cells.Item(48, 4).Value = "'0.1233"
$ws.Cells.Item(48, 4).Style = 'Normal'
$ws.Cells.Item(48, 5).Value = '  +1.37%  '

$ws.Cells.Item(49, 4).Value = "'0.05815"
$ws.Cells.Item(49, 4).Style = 'Normal'
$ws.Cells.Item(49, 5).Value = '  -1.27%  '

$ws.Cells.Item(50, 4).Value = "'34.70"
$ws.Cells.Item(50, 4).Style = 'Normal'
$ws.Cells.Item(50, 5).Value = '  -1.05%  '

$ws.Cells.Item(51, 4).Value = "'0.8945"
$ws.Cells.Item(51, 4).Style = 'Normal'
$ws.Cells.Item(51, 5).Value = '  +0.91%  '
